$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nuit")

# "Catégorie" column (P) on the "Nuit" sheet is being re-classified: a new
# "Multi-usage" category is introduced (it lands in the shared-strings table
# right after "CVC", which is what shifts "Eclairage" further down the table
# and changes P5's resolved text from "Eclairage" to "Multi-usage" even
# though its stored string index doesn't move), and several rows are
# reassigned to "Investigation en cours" / "Eclairage".
$ws.Range("P5").Value  = "Multi-usage"
$ws.Range("P6").Value  = "Investigation en cours"
$ws.Range("P7").Value  = "Eclairage"
$ws.Range("P10").Value = "Investigation en cours"
$ws.Range("P11").Value = "Investigation en cours"
$ws.Range("P14").Value = "Eclairage"
